$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation (2023-04-11) was added at the top of the
# data (row 224), pushing all subsequent rows down by one. The new row
# carries the same price figures as the most recent existing entry
# (previously the last row of the sheet).
$ws.Rows.Item(224).EntireRow.Insert()

$ws.Range("A224").Value = 6
$ws.Range("B224").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C224").Value = "Metropolitana"
$ws.Range("D224").Value = "2023-04-11"
$ws.Range("E224").Value = 13
$ws.Range("F224").Value = 100112029
$ws.Range("G224").Value = "Orégano"
$ws.Range("H224").Value = "Sin especificar"
$ws.Range("I224").Value = "Primera"
$ws.Range("J224").Value = 46
$ws.Range("K224").Value = 17000
$ws.Range("L224").Value = 18000
$ws.Range("M224").Value = 17478
$ws.Range("N224").Value = "$/docena de atados"
$ws.Range("O224").Value = "Región Metropolitana"
$ws.Range("P224").Value = 5826
$ws.Range("Q224").Value = 3
$ws.Range("R224").Value = "Hortaliza"
